$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990868838/1749474758173/gFRHQATQ91-camera-video-segment-driver-1749474760673.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSIF64JU42%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T200022Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEOr%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDxIePx%2FsJ6%2FbMwqnOPO9B6Hj%2B0VLeRVtXvmy9Kuq0fiQIgJCQHNAIOGGfhMvoVuXdcfwqKETJ5LQmn5BbJqKEtYqsq5gMIw%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDJ5bb6ZUGpNiTPgzOCq6A5%2FpFb1CnQwrHeedyB0DYvSwxyLfh0mtOfGdoJibCjAhqoxFjEGRbXZtZkYBIvCiGNQFAPSZ29jor2e%2BZCCuzfOyfV0xT4AU5wWyWesKG44aKH8eRLU5G8rQsqV9yXSmk9%2Bm%2FrJsNoxB%2B4aO8NGy3YxhQQLoF%2FLa6pNmQ9Hi%2FMBvzRFMUQvJmZZpX3PNiIH%2B4tEXxjAHLaLLft6%2FpxLHqwWN25%2FObE6aimZNHxomYH%2BlCcqE%2BGbsZrYsvbmP%2BYPjw2n4iCm7KCc%2BV0%2F6oqicHQPZwHatQiWn8DxbCPjhROxItBVOg3Etm3rmZa8EawXAj%2BVSGsknS4ISjx7X2pncaOgJDeQdps9xwl%2B6o68794sXBQ6UXqTSz14K6%2FaJ7EUzW9d644bYQMDGTzDOUp4egdMaJm5RtgcrdwLahS98XSW89FYdh%2BAuYXnclWZ6cI2BuZU7uletT2t7vH9oXHisFqUVC5PTAvzIYi%2FDnrkUcTNSBOpOL1n9NO47IoMvMS%2BtIAtKGEKPixXatNaAf7%2BndYgkApgJ3oMMwfp84nAMTwuIGSigfl4HZsRZfjxD%2Frh12xfc%2FFxtojKhGRgwzOihwgY6pQHlMiuzK4pZKZ3saBmKZ8sO78l7Tvg5RDYMYF7gzjgxK6YSqujcpiny0F0vMBVqAb5U1n1lHOI30TyP04YO1Sb2vzUWVWVud8B9OF8ojJMBjzdsM0%2FVDCUVeemOfOEOeRwX5kXvxqyIHEdFhtMpyc2hDegiTFzoanB7rjhJPu4DHejnLbqAf95w03dNc3Bp6RIU12XmamqFtC%2Fb6Y9euwWDZ9PzZRc%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2011%20Jun%202025%2004%3A00%3A22%20GMT&X-Amz-Signature=fae91a24bf121fdc3171099939a34b1479c2350da437eb3a3079fff0df747b6c'
$ws.Range("L2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990868838/1749474758173/7Z3ova6Wtr-camera-video-segment-1749474760673.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSIF64JU42%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T200022Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEOr%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDxIePx%2FsJ6%2FbMwqnOPO9B6Hj%2B0VLeRVtXvmy9Kuq0fiQIgJCQHNAIOGGfhMvoVuXdcfwqKETJ5LQmn5BbJqKEtYqsq5gMIw%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDJ5bb6ZUGpNiTPgzOCq6A5%2FpFb1CnQwrHeedyB0DYvSwxyLfh0mtOfGdoJibCjAhqoxFjEGRbXZtZkYBIvCiGNQFAPSZ29jor2e%2BZCCuzfOyfV0xT4AU5wWyWesKG44aKH8eRLU5G8rQsqV9yXSmk9%2Bm%2FrJsNoxB%2B4aO8NGy3YxhQQLoF%2FLa6pNmQ9Hi%2FMBvzRFMUQvJmZZpX3PNiIH%2B4tEXxjAHLaLLft6%2FpxLHqwWN25%2FObE6aimZNHxomYH%2BlCcqE%2BGbsZrYsvbmP%2BYPjw2n4iCm7KCc%2BV0%2F6oqicHQPZwHatQiWn8DxbCPjhROxItBVOg3Etm3rmZa8EawXAj%2BVSGsknS4ISjx7X2pncaOgJDeQdps9xwl%2B6o68794sXBQ6UXqTSz14K6%2FaJ7EUzW9d644bYQMDGTzDOUp4egdMaJm5RtgcrdwLahS98XSW89FYdh%2BAuYXnclWZ6cI2BuZU7uletT2t7vH9oXHisFqUVC5PTAvzIYi%2FDnrkUcTNSBOpOL1n9NO47IoMvMS%2BtIAtKGEKPixXatNaAf7%2BndYgkApgJ3oMMwfp84nAMTwuIGSigfl4HZsRZfjxD%2Frh12xfc%2FFxtojKhGRgwzOihwgY6pQHlMiuzK4pZKZ3saBmKZ8sO78l7Tvg5RDYMYF7gzjgxK6YSqujcpiny0F0vMBVqAb5U1n1lHOI30TyP04YO1Sb2vzUWVWVud8B9OF8ojJMBjzdsM0%2FVDCUVeemOfOEOeRwX5kXvxqyIHEdFhtMpyc2hDegiTFzoanB7rjhJPu4DHejnLbqAf95w03dNc3Bp6RIU12XmamqFtC%2Fb6Y9euwWDZ9PzZRc%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2011%20Jun%202025%2004%3A00%3A22%20GMT&X-Amz-Signature=4e222ba32f9c2d09f67aad94776b23a45cce335042afa96730b6dd30a57d50ef'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474990868838/1749474688613/JLWgs3viB6-camera-video-segment-driver-1749474691113.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSIF64JU42%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T200022Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEOr%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDxIePx%2FsJ6%2FbMwqnOPO9B6Hj%2B0VLeRVtXvmy9Kuq0fiQIgJCQHNAIOGGfhMvoVuXdcfwqKETJ5LQmn5BbJqKEtYqsq5gMIw%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDJ5bb6ZUGpNiTPgzOCq6A5%2FpFb1CnQwrHeedyB0DYvSwxyLfh0mtOfGdoJibCjAhqoxFjEGRbXZtZkYBIvCiGNQFAPSZ29jor2e%2BZCCuzfOyfV0xT4AU5wWyWesKG44aKH8eRLU5G8rQsqV9yXSmk9%2Bm%2FrJsNoxB%2B4aO8NGy3YxhQQLoF%2FLa6pNmQ9Hi%2FMBvzRFMUQvJmZZpX3PNiIH%2B4tEXxjAHLaLLft6%2FpxLHqwWN25%2FObE6aimZNHxomYH%2BlCcqE%2BGbsZrYsvbmP%2BYPjw2n4iCm7KCc%2BV0%2F6oqicHQPZwHatQiWn8DxbCPjhROxItBVOg3Etm3rmZa8EawXAj%2BVSGsknS4ISjx7X2pncaOgJDeQdps9xwl%2B6o68794sXBQ6UXqTSz14K6%2FaJ7EUzW9d644bYQMDGTzDOUp4egdMaJm5RtgcrdwLahS98XSW89FYdh%2BAuYXnclWZ6cI2BuZU7uletT2t7vH9oXHisFqUVC5PTAvzIYi%2FDnrkUcTNSBOpOL1n9NO47IoMvMS%2BtIAtKGEKPixXatNaAf7%2BndYgkApgJ3oMMwfp84nAMTwuIGSigfl4HZsRZfjxD%2Frh12xfc%2FFxtojKhGRgwzOihwgY6pQHlMiuzK4pZKZ3saBmKZ8sO78l7Tvg5RDYMYF7gzjgxK6YSqujcpiny0F0vMBVqAb5U1n1lHOI30TyP04YO1Sb2vzUWVWVud8B9OF8ojJMBjzdsM0%2FVDCUVeemOfOEOeRwX5kXvxqyIHEdFhtMpyc2hDegiTFzoanB7rjhJPu4DHejnLbqAf95w03dNc3Bp6RIU12XmamqFtC%2Fb6Y9euwWDZ9PzZRc%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2011%20Jun%202025%2004%3A00%3A22%20GMT&X-Amz-Signature=20354272f32eaf97a593b8543f4e071f6cf75253066eb19c46d518e245043ac8'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474990868838/1749474688613/LOHXs2uAgt-camera-video-segment-1749474691113.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSIF64JU42%2F20250610%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250610T200022Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEOr%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDxIePx%2FsJ6%2FbMwqnOPO9B6Hj%2B0VLeRVtXvmy9Kuq0fiQIgJCQHNAIOGGfhMvoVuXdcfwqKETJ5LQmn5BbJqKEtYqsq5gMIw%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDJ5bb6ZUGpNiTPgzOCq6A5%2FpFb1CnQwrHeedyB0DYvSwxyLfh0mtOfGdoJibCjAhqoxFjEGRbXZtZkYBIvCiGNQFAPSZ29jor2e%2BZCCuzfOyfV0xT4AU5wWyWesKG44aKH8eRLU5G8rQsqV9yXSmk9%2Bm%2FrJsNoxB%2B4aO8NGy3YxhQQLoF%2FLa6pNmQ9Hi%2FMBvzRFMUQvJmZZpX3PNiIH%2B4tEXxjAHLaLLft6%2FpxLHqwWN25%2FObE6aimZNHxomYH%2BlCcqE%2BGbsZrYsvbmP%2BYPjw2n4iCm7KCc%2BV0%2F6oqicHQPZwHatQiWn8DxbCPjhROxItBVOg3Etm3rmZa8EawXAj%2BVSGsknS4ISjx7X2pncaOgJDeQdps9xwl%2B6o68794sXBQ6UXqTSz14K6%2FaJ7EUzW9d644bYQMDGTzDOUp4egdMaJm5RtgcrdwLahS98XSW89FYdh%2BAuYXnclWZ6cI2BuZU7uletT2t7vH9oXHisFqUVC5PTAvzIYi%2FDnrkUcTNSBOpOL1n9NO47IoMvMS%2BtIAtKGEKPixXatNaAf7%2BndYgkApgJ3oMMwfp84nAMTwuIGSigfl4HZsRZfjxD%2Frh12xfc%2FFxtojKhGRgwzOihwgY6pQHlMiuzK4pZKZ3saBmKZ8sO78l7Tvg5RDYMYF7gzjgxK6YSqujcpiny0F0vMBVqAb5U1n1lHOI30TyP04YO1Sb2vzUWVWVud8B9OF8ojJMBjzdsM0%2FVDCUVeemOfOEOeRwX5kXvxqyIHEdFhtMpyc2hDegiTFzoanB7rjhJPu4DHejnLbqAf95w03dNc3Bp6RIU12XmamqFtC%2Fb6Y9euwWDZ9PzZRc%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2011%20Jun%202025%2004%3A00%3A22%20GMT&X-Amz-Signature=118a46cb4faf7da524e3942a5b548b9b0a01f9d681a74285fa01de681d3c8059'
